# Append two new order lines (rows 21-22) to the Bakery order sheet.
# The existing sheet stores every value (including numeric-looking
# SKU/Quantity/Cost figures) as literal text, so format the target
# cells as Text *before* writing the values - this mirrors typing into
# a Text-formatted cell in Excel and keeps "38505", "1", "56.85", etc.
# from being auto-converted into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21:E22").NumberFormat = "@"

$ws.Range("A21").Value = "38505"
$ws.Range("B21").Value = "Container - Paper Clamshell (Avocado Box)"
$ws.Range("C21").Value = "1"
$ws.Range("D21").Value = "56.85"
$ws.Range("E21").Value = "56.85"

$ws.Range("A22").Value = "22517"
$ws.Range("B22").Value = "Bag Sheet Pan Cover 30x43"
$ws.Range("C22").Value = "2"
$ws.Range("D22").Value = "27.77"
$ws.Range("E22").Value = "55.54"
